# Update of seepage model (EVT)
# Inserts two new log entries at the top of the Stream_seepage sheet's
# worklog table (rows 2 and 3), pushing the existing entries down by two
# rows, and updates the sheet selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stream_seepage")

# Make room for the two new entries right under the header row.
$ws.Rows("2:3").Insert()

# --- New row 2: "Set up" entry (2023-04-26) ---
$ws.Range("A2").Value = 45042
$ws.Range("A5").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("B2").Value = "Set up"
$ws.Range("C2").Value = "I removed the steady state period to avoid excessive drawdown on start up with new EVT input. No issues with water budget, 1 or two random steps with 2% error. Good fit with hydrographs (NSE>0.5). The difference again though is that there is no baseflow occuring for most time steps which seems to be the fault of excessive ET or perhaps the fault of a change in stream elevation with the new model grid. I'm inclined to think it's an issue with the model grid being 2m steps now so the stream bottom is potentially being weird, but again it could just be higher ET rates with that 10 m rooting depth next to the stream."
$ws.Range("D2").Value = "One thing that might be worth addressing is whether we are overestimating ET by not including a lot pumping, this is only a concern if ET is evaluated in the results."
$ws.Rows("2").RowHeight = 86.4

# --- New row 3: "EVT" entry (2023-04-26) ---
$ws.Range("A3").Value = 45042
$ws.Range("A5").Copy()
$ws.Range("A3").PasteSpecial(-4122)

$ws.Range("B3").Value = "EVT"
$ws.Range("C3").Value = "I removed EVT under the stream to try to improve the steady state and avoid excess drawdown directly below the stream but this led to the weird effect again where the steady state levels start really low. But removing EVT below SFR also improved runtime from 1 hr 20 to 28 min."
$ws.Rows("3").RowHeight = 43.2

$excel.CutCopyMode = 0

# Matches the cursor position left behind in the edited workbook.
$ws.Range("D3").Select() | Out-Null
